$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values, preserving them as text
# (matching the original inline-string cell representation in the workbook).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '329.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '4.86%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.37'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.16%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.757'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '12.69%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08116'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.44%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.592'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.84%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.782'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.37%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.973'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4.68%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.12%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9459'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.46%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1310'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.27%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.020'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '39.97%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09272'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.43%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.03444'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '4.40%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.09618'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.70%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001312'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-5.76%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006330'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '5.28%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.360'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.64%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3563'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.08%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '9.41%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2412'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '5.19%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04434'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.16%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001261'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '5.77%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004368'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.25%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001092'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-17.36%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003995'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '1.12%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02465'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '8.22%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05307'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.48%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007480'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.23%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1434'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.18%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008857'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.30%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002064'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.90%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01041'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '31.43%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006887'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '9.11%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.78%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003513'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '23.34%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001702'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '1.51%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.78%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.78%'
